$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New fit result row: Fit012 (inserted before Fit016, appended at the end of the table)
$ws.Range("A18").Value = "Fit012"
$ws.Range("B18").Value = "D2"
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = "ERF"
$ws.Range("F18").Value = "None"
$ws.Range("G18").Value = 1.372534128271123
$ws.Range("H18").Value = 0.03527432770636527
$ws.Range("I18").Value = 0.3165515876592619
$ws.Range("J18").Value = 0.009680758460240567
$ws.Range("K18").Value = 7
$ws.Range("L18").Value = 10
$ws.Range("M18").Value = 0.009680758460240567
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0.0000000000003449748154796696

# New fit result row: Fit016 (exp energy based cutoff function fit)
$ws.Range("A19").Value = "Fit016"
$ws.Range("B19").Value = "D2"
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = "ERF"
$ws.Range("F19").Value = "None"
$ws.Range("G19").Value = 1.633497756077309
$ws.Range("H19").Value = 0.05397988103115311
$ws.Range("I19").Value = 0.3573013357182082
$ws.Range("J19").Value = 0.01576847984847383
$ws.Range("K19").Value = 7
$ws.Range("L19").Value = 10
$ws.Range("M19").Value = 0.01576847984847383
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 0.000000000003044860279865397
